# Build site at 2022-09-26 16:07:08 UTC
# This script reproduces the edit made to LOM3110.xlsx:
#  - the "Objetivos:" value (row 10, B/C) is replaced with a professor id/name string
#  - the three standalone rows holding only the professor id/name values
#    (old rows 13-15, which had no label in column A) are removed, shifting
#    everything below them up by three rows
#  - after the shift, the B/C "value" cell for several of the remaining
#    labelled rows (Programa resumido, Programa, Método, Critério,
#    Norma de recuperação, Bibliografia) is updated to hold the text that
#    now lines up with that label

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Row 10 ("Objetivos:") no longer shows the long objectives paragraph;
#    it now shows the first professor's id/name.
$ws.Range("B10").Value = "7459752 - Maria Ismenia Sodero Toledo Faria"
$ws.Range("C10").Value = "7459752 - Maria Ismenia Sodero Toledo Faria"

# 2) Remove the three bare professor-name rows (old rows 13, 14, 15).
#    Deleting the same row index three times in a row removes all three,
#    because each deletion shifts the following rows up into that index.
$ws.Rows.Item(13).Delete() | Out-Null
$ws.Rows.Item(13).Delete() | Out-Null
$ws.Rows.Item(13).Delete() | Out-Null

# 3) Fix up the B/C values on the rows that shifted into the old
#    "Programa resumido:"/"Programa:"/"Método:"/"Critério:"/
#    "Norma de recuperação:"/"Bibliografia:" positions (now rows 13, 15,
#    18, 19, 20, 21) so each label is paired with the correct value.
$ws.Range("B13").Value = "7459752 - Maria Ismenia Sodero Toledo Faria"
$ws.Range("C13").Value = "7459752 - Maria Ismenia Sodero Toledo Faria"

$ws.Range("B15").Value = "2166002 - Sandra Giacomin Schneider"
$ws.Range("C15").Value = "2166002 - Sandra Giacomin Schneider"

$ws.Range("B18").Value = "1922320 - Sebastiao Ribeiro"
$ws.Range("C18").Value = "1922320 - Sebastiao Ribeiro"

$metodoText = 'O método utilizado tem por fundamento a aprendizagem baseada em projetos que visa desenvolver as competências técnicas relativas ao tema do projeto, bem como competências transversais, tais como: aprender a aprender, trabalho em equipe, relacionamento interpessoal, capacidade de comunicação oral e verbal e aspectos de liderança, dentre outros.Os alunos serão divididos em grupos que desenvolverão um projeto durante o semestre relacionado a um tema de Engenharia de Materiais, similar ao que eles irão encontrar na vida real no efetivo exercício de sua profissão.Cada grupo deverá buscar o conhecimento prático necessário para ser aplicado no desenvolvimento do projeto.As aulas ocorrerão por meio de uma reunião da equipe de trabalho para tratar do projeto; palestras e dinâmicas relativas ao tema do projeto, conduzidas por professores ou profissionais de empresas.'
$ws.Range("B19").Value = $metodoText
$ws.Range("C19").Value = $metodoText

$notaText = 'A nota será individual e será a média ponderada de entregas do projeto, tais como: projeto preliminar, projeto final, envolvimento do aluno com o projeto, avaliação dos pares, autoavaliação e apresentação de trabalhos, dentre outros.O detalhamento dos pesos para ponderação da média da disciplina será definido por uma equipe de professores que atuarão na avaliação da disciplina.'
$ws.Range("B20").Value = $notaText
$ws.Range("C20").Value = $notaText

$ws.Range("B21").Value = "não há"
$ws.Range("C21").Value = "não há"
